$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the second paragraph entirely:
#    "This is the document you will need to change.  Delete everything
#     below the above line."
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(2)
if ($introPara.Range.Text -like "This is the document you will need to change.*") {
    $introPara.Range.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 2. The paragraph that used to be the third paragraph is now the
#    second paragraph. Across several runs it currently reads:
#      "Add instructions for your tutor to pull (merge) your request
#       to the mainline.  Note that the tutor will not merge all
#       requests, but you should say how it could be achiev[_GoBack]ed."
#    Replace that whole sentence with the new text, keeping the
#    "_GoBack" bookmark exactly where it was (immediately after the
#    remaining text, just as it previously sat right before "ed.").
# ------------------------------------------------------------------
$p = $d.Paragraphs(2)
$pStart = $p.Range.Start
$pEnd = $p.Range.End

$bm = $d.Bookmarks("_GoBack")

# Delete any leftover text after the bookmark (originally "ed."),
# keeping the paragraph mark itself intact.
$afterBookmark = $d.Range($bm.End, $pEnd - 1)
if ($afterBookmark.Start -lt $afterBookmark.End) {
    $afterBookmark.Delete() | Out-Null
}

# Replace all of the text before the bookmark with the new sentence.
$beforeBookmark = $d.Range($pStart, $bm.Start)
$beforeBookmark.Text = "Changes to this document have been made from bgreen21 account"
